# Update cell values per commit: "update scripts wuth new tpm"
# Re-run of the NATMI TPM pipeline produced refreshed ligand/receptor
# expression statistics for the Fzd8-Ckap4 sheet (rows 2-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.083576666666666
$ws.Range("H2").Value = 9.250729999999999
$ws.Range("I2").Value = 0.2272509363535097
$ws.Range("J2").Value = 0.2272509363535097
$ws.Range("M2").Value = 1.695728
$ws.Range("N2").Value = 5.087184000000001
$ws.Range("O2").Value = 0.06675079911082282
$ws.Range("P2").Value = 0.06675079911082282
$ws.Range("Q2").Value = 5.228907293813333
$ws.Range("R2").Value = 47.06016564432
$ws.Range("S2").Value = 0.01516918160027951
$ws.Range("T2").Value = 0.01516918160027951

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.083576666666666
$ws.Range("H3").Value = 9.250729999999999
$ws.Range("I3").Value = 0.2272509363535097
$ws.Range("J3").Value = 0.2272509363535097
$ws.Range("O3").Value = 0.7828150713365326
$ws.Range("P3").Value = 0.7828150713365327
$ws.Range("Q3").Value = 61.32162447108333
$ws.Range("R3").Value = 551.89462023975
$ws.Range("S3").Value = 0.1778954579528665
$ws.Range("T3").Value = 0.1778954579528665

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.083576666666666
$ws.Range("H4").Value = 9.250729999999999
$ws.Range("I4").Value = 0.2272509363535097
$ws.Range("J4").Value = 0.2272509363535097
$ws.Range("M4").Value = 3.821607666666667
$ws.Range("N4").Value = 11.464823
$ws.Range("O4").Value = 0.1504341295526447
$ws.Range("P4").Value = 0.1504341295526447
$ws.Range("Q4").Value = 11.78422023008778
$ws.Range("R4").Value = 106.05798207079
$ws.Range("S4").Value = 0.03418629680036368
$ws.Range("T4").Value = 0.03418629680036368

# Row 5
$ws.Range("G5").Value = 6.453984666666667
$ws.Range("I5").Value = 0.4756405360586227
$ws.Range("J5").Value = 0.4756405360586227
$ws.Range("M5").Value = 1.695728
$ws.Range("N5").Value = 5.087184000000001
$ws.Range("O5").Value = 0.06675079911082282
$ws.Range("P5").Value = 0.06675079911082282
$ws.Range("Q5").Value = 10.94420251083733
$ws.Range("R5").Value = 98.49782259753601
$ws.Range("S5").Value = 0.0317493858714132
$ws.Range("T5").Value = 0.0317493858714132

# Row 6
$ws.Range("G6").Value = 6.453984666666667
$ws.Range("I6").Value = 0.4756405360586227
$ws.Range("J6").Value = 0.4756405360586227
$ws.Range("O6").Value = 0.7828150713365326
$ws.Range("P6").Value = 0.7828150713365327
$ws.Range("Q6").Value = 128.3473274232834
$ws.Range("S6").Value = 0.3723385801652774
$ws.Range("T6").Value = 0.3723385801652774

# Row 7
$ws.Range("G7").Value = 6.453984666666667
$ws.Range("I7").Value = 0.4756405360586227
$ws.Range("J7").Value = 0.4756405360586227
$ws.Range("M7").Value = 3.821607666666667
$ws.Range("N7").Value = 11.464823
$ws.Range("O7").Value = 0.1504341295526447
$ws.Range("P7").Value = 0.1504341295526447
$ws.Range("Q7").Value = 24.66459728268245
$ws.Range("R7").Value = 221.981375544142
$ws.Range("S7").Value = 0.0715525700219322
$ws.Range("T7").Value = 0.07155257002193219

# Row 8
$ws.Range("G8").Value = 4.031477000000001
$ws.Range("H8").Value = 12.094431
$ws.Range("I8").Value = 0.2971085275878677
$ws.Range("J8").Value = 0.2971085275878677
$ws.Range("M8").Value = 1.695728
$ws.Range("N8").Value = 5.087184000000001
$ws.Range("O8").Value = 0.06675079911082282
$ws.Range("P8").Value = 0.06675079911082282
$ws.Range("Q8").Value = 6.836288430256001
$ws.Range("R8").Value = 61.52659587230401
$ws.Range("S8").Value = 0.01983223163913012
$ws.Range("T8").Value = 0.01983223163913012

# Row 9
$ws.Range("G9").Value = 4.031477000000001
$ws.Range("H9").Value = 12.094431
$ws.Range("I9").Value = 0.2971085275878677
$ws.Range("J9").Value = 0.2971085275878677
$ws.Range("O9").Value = 0.7828150713365326
$ws.Range("P9").Value = 0.7828150713365327
$ws.Range("Q9").Value = 80.17206814742502
$ws.Range("R9").Value = 721.5486133268251
$ws.Range("S9").Value = 0.2325810332183888
$ws.Range("T9").Value = 0.2325810332183888

# Row 10
$ws.Range("G10").Value = 4.031477000000001
$ws.Range("H10").Value = 12.094431
$ws.Range("I10").Value = 0.2971085275878677
$ws.Range("J10").Value = 0.2971085275878677
$ws.Range("M10").Value = 3.821607666666667
$ws.Range("N10").Value = 11.464823
$ws.Range("O10").Value = 0.1504341295526447
$ws.Range("P10").Value = 0.1504341295526447
$ws.Range("Q10").Value = 15.40672341119034
$ws.Range("R10").Value = 138.660510700713
$ws.Range("S10").Value = 0.04469526273034878
$ws.Range("T10").Value = 0.04469526273034878
